$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing A:E data to B:F.
$ws.Columns.Item(1).Insert()

# New header for column A.
$ws.Cells.Item(1,1).Value = "ID"

# Give the new A1 header cell the same formatting (bold, bordered,
# centered) as the other header cells by copying B1's format onto it.
$ws.Cells.Item(1,2).Copy()
$ws.Cells.Item(1,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New ID values for rows 2-16.
$ids = @("Hb 47", "Hb 2", "Hb 3", "S 24", "S 28", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 3", "S 4", "S 5", "Hb 74", "Hb 79")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
